# correção nos dados e inicio da analise PNAD 2009
#
# The original "unnamed: 1_level_1" / "unnamed: 5_level_1" placeholder
# labels (leftovers from a pandas multi-index header) are replaced by
# "total" in the two spots that used them, which also makes those two
# shared strings unused so they drop out of the table on save.
# The row-4..row-20 label column was off by one category (each label had
# been paired with the wrong region name) — shift each A-column label to
# the correct adjoining region name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sub-header row. B2/C2/F2 collapse to "total"; D2/E2 keep their
# existing text (condição de ocupação na semana de referência[.1]).
$ws.Cells.Item(2, 2).Value = "total"
$ws.Cells.Item(2, 3).Value = "total"
$ws.Cells.Item(2, 4).Value = "condição de ocupação na semana de referência"
$ws.Cells.Item(2, 5).Value = "condição de ocupação na semana de referência.1"
$ws.Cells.Item(2, 6).Value = "total"

# Column A labels (rows 4-20): corrected to the right region names.
$ws.Cells.Item(4, 1).Value = "pará"
$ws.Cells.Item(5, 1).Value = "região metropolitana de belém"
$ws.Cells.Item(6, 1).Value = "ceará"
$ws.Cells.Item(7, 1).Value = "região metropolitana de fortaleza"
$ws.Cells.Item(8, 1).Value = "pernambuco"
$ws.Cells.Item(9, 1).Value = "região metropolitana de recife"
$ws.Cells.Item(10, 1).Value = "bahia"
$ws.Cells.Item(11, 1).Value = "região metropolitana de salvador"
$ws.Cells.Item(12, 1).Value = "minas gerais"
$ws.Cells.Item(13, 1).Value = "região metropolitana de belo horizonte"
$ws.Cells.Item(14, 1).Value = "      rio de janeiro"
$ws.Cells.Item(15, 1).Value = "região metropolitana do rio de janeiro"
$ws.Cells.Item(16, 1).Value = "são paulo"
$ws.Cells.Item(17, 1).Value = "região metropolitana de são paulo"
$ws.Cells.Item(18, 1).Value = "       paraná"
$ws.Cells.Item(19, 1).Value = "  região metropolitana de curitiba"
$ws.Cells.Item(20, 1).Value = "      rio grande do sul"
